# Changed Check Date and Service date for N70
#
# Input.xlsx — B2 held a numeric invoice number (191550000307). It is
# replaced with a new invoice number (191480000297) stored as *text*
# (so it round-trips through sharedStrings.xml as a shared string
# rather than a numeric <v>), and the cell gets a distinct look:
# a 10pt font, a thin box border on all four sides, and vertically
# centered text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B2")

# Apply the new formatting first (font size, vertical centering, thin
# border on all sides) while the cell still holds its old numeric
# value/format — this mutates the single style record B2 already owns
# instead of allocating a brand new one, which matches how the
# existing xf entry (s="2") was rewritten in place rather than a new
# style slot being appended.
$cell.Font.Size = 10
$cell.VerticalAlignment = -4108   # xlCenter
$cell.Borders.LineStyle = 1       # xlContinuous -> thin box border

# Now force the cell to Text format and write the new invoice number
# as a text string (keeps leading context intact and forces it into
# sharedStrings.xml rather than being stored as a <v> number).
$cell.NumberFormat = "@"
$cell.Value = "191480000297"
